{"js": "// Update the division \"fact family\" answers in the first table.\n// The replacements are positional (several old/new values repeat),\n// so we walk the table's cells in document order and overwrite the\n// Nth non-empty cell's text with the Nth replacement value. Using\n// table.values (get + set) preserves each cell's existing paragraph/\n// run formatting and only swaps the <w:t> text, matching the diff.\n\nconst replacements = [\n  \"31\u00f73=10, 1\", \"50\u00f78=6, 2\", \"13\u00f72=6, 1\", \"21\u00f72=10, 1\", \"18\u00f78=2, 2\",\n  \"27\u00f78=3, 3\", \"89\u00f77=12, 5\", \"21\u00f76=3, 3\", \"49\u00f75=9, 4\", \"69\u00f75=13, 4\",\n  \"62\u00f76=10, 2\", \"26\u00f78=3, 2\", \"99\u00f74=24, 3\", \"83\u00f76=13, 5\", \"72\u00f76=12, 0\",\n  \"53\u00f77=7, 4\", \"89\u00f72=44, 1\", \"13\u00f72=6, 1\", \"67\u00f79=7, 4\", \"44\u00f73=14, 2\",\n  \"34\u00f72=17, 0\", \"82\u00f77=11, 5\", \"18\u00f78=2, 2\", \"22\u00f72=11, 0\", \"22\u00f74=5, 2\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values.map((row) => row.slice());\n\nlet k = 0;\nfor (let r = 0; r < values.length && k < replacements.length; r++) {\n  for (let c = 0; c < values[r].length && k < replacements.length; c++) {\n    if (values[r][c] !== \"\") {\n      values[r][c] = replacements[k];\n      k++;\n    }\n  }\n}\n\ntable.values = values;\nawait context.sync();\n", "ps1": "# Update the division \"fact family\" answers in the first table.\n# The replacements are positional (several old/new values repeat),\n# so we walk the table's cells in document order (row-major) and\n# overwrite the Nth non-empty cell's text with the Nth replacement\n# value. Assigning Cell.Range.Text preserves the cell's existing\n# paragraph/run formatting and only swaps the <w:t> text, matching\n# the diff.\n\n$replacements = @(\n  \"31\u00f73=10, 1\", \"50\u00f78=6, 2\", \"13\u00f72=6, 1\", \"21\u00f72=10, 1\", \"18\u00f78=2, 2\",\n  \"27\u00f78=3, 3\", \"89\u00f77=12, 5\", \"21\u00f76=3, 3\", \"49\u00f75=9, 4\", \"69\u00f75=13, 4\",\n  \"62\u00f76=10, 2\", \"26\u00f78=3, 2\", \"99\u00f74=24, 3\", \"83\u00f76=13, 5\", \"72\u00f76=12, 0\",\n  \"53\u00f77=7, 4\", \"89\u00f72=44, 1\", \"13\u00f72=6, 1\", \"67\u00f79=7, 4\", \"44\u00f73=14, 2\",\n  \"34\u00f72=17, 0\", \"82\u00f77=11, 5\", \"18\u00f78=2, 2\", \"22\u00f72=11, 0\", \"22\u00f74=5, 2\"\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$k = 0\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n  for ($c = 1; $c -le $table.Columns.Count; $c++) {\n    if ($k -ge $replacements.Length) { break }\n    $cell = $table.Cell($r, $c)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne \"\") {\n      $cell.Range.Text = $replacements[$k]\n      $k = $k + 1\n    }\n  }\n}\n"}
